# Chess - V 1.5 Fixed Box Issue
#
# The design.docx canvas gained a new "Rectangle: Rounded Corners 1" callout box
# (an mc:AlternateContent drawing with a DrawingML roundRect + VML v:roundrect
# fallback) inserted as a brand-new run right before the very first picture run.
# Word also (re-)stamped *every* <wp:anchor> already in the document with a fresh
# wp14:anchorId/wp14:editId pair (these are opaque undo/redo correlation ids that
# Word mints whenever it resaves a drawing host paragraph). We reproduce both
# effects by editing the raw package XML through Range.WordOpenXML, which gives us
# byte-for-byte control over the OOXML (the Word object model has no direct
# "add rounded-rectangle AutoShape with this exact XML" verb).

$d = $word.ActiveDocument
$rng = $d.Content
$xml = $rng.WordOpenXML

# --- 1) Insert the new rounded-rectangle shape run just before the first picture
#        run, and stamp that first picture run anchor with its new ids. ---
$oldFirstRunStart = @'
<w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251640822" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$newFirstRunStart = @'
<w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251711488" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="6D306902" wp14:editId="3822DDCA"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>left</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>5009712</wp:posOffset></wp:positionV><wp:extent cx="1036101" cy="1045968"/><wp:effectExtent l="19050" t="19050" r="12065" b="20955"/><wp:wrapNone/><wp:docPr id="684123505" name="Rectangle: Rounded Corners 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1036101" cy="1045968"/></a:xfrm><a:prstGeom prst="roundRect"><a:avLst><a:gd name="adj" fmla="val 3630"/></a:avLst></a:prstGeom><a:noFill/><a:ln w="38100" cap="flat" cmpd="sng" algn="ctr"><a:solidFill><a:schemeClr val="accent6"/></a:solidFill><a:prstDash val="solid"/><a:round/><a:headEnd type="none" w="med" len="med"/><a:tailEnd type="none" w="med" len="med"/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:lnRef><a:fillRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:fillRef><a:effectRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="accent6"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:roundrect w14:anchorId="15C0A637" id="Rectangle: Rounded Corners 1" o:spid="_x0000_s1026" style="position:absolute;margin-left:0;margin-top:394.45pt;width:81.6pt;height:82.35pt;z-index:251711488;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:left;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:middle" arcsize="2379f" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQA9RknDyQIAAB8GAAAOAAAAZHJzL2Uyb0RvYy54bWysVNtOGzEQfa/Uf7D8XrJLQgoRGxRBqSqh&#10;goCKZ8drZ7fyrbZz69f32LtJoBepVH3xjncunjlnZs4vNlqRlfChtaai5VFBiTDc1q1ZVPTL4/W7&#10;U0pCZKZmyhpR0a0I9GL69s352k3EsW2sqoUnCGLCZO0q2sToJoNB4I3QLBxZJwyU0nrNIq5+Mag9&#10;WyO6VoPjohgP1tbXzlsuQsDfq05Jpzm+lILHWymDiERVFLnFfPp8ztM5mJ6zycIz17S8T4P9Qxaa&#10;tQaP7kNdscjI0re/hNIt9zZYGY+41QMrZctFrgHVlMVP1Tw0zIlcC8AJbg9T+H9h+efVg7vzgGHt&#10;wiRATFVspNfpi/zIJoO13YMlNpFw/CyL4bgsSko4dGUxOjkbnyY4Bwd350P8KKwmSaiot0tT34OS&#10;jBRb3YSYIauJYRq9weqvlEitQMCKKTIcDzM/CNjbQtqFTI7GXrdKZQaVIeuKDk/LAiRzhkaSikWI&#10;2tUVDWZBCVMLdCiPPr8erGrr5J4C5W4Tl8oTvIs8OBcmjvtiXlim569YaDrDrOpaKNeWU2kEqz+Y&#10;msStQ00GXU9TblrUlCiBFJKULSNr1d9YomxlAOyBoSzFrRIpe2XuhSRtnYnqyvGLeaqm63MMIkDZ&#10;dXsOBodkKFH/K317l+Qt8ni90n/vlN+3Ju79dWtsz00a/j/RITufHRwdCAmPua23d6jZdjMeHL9u&#10;wdYNC/GOefQUQMCiirc4pLKgxPYSJY3133/3P9lj1qAFh1gSaKVvS+bBqPpkMIVn5WiUtkq+jE7e&#10;H+Pin2vmzzVmqS8t+gszg+yymOyj2onSW/2EfTZLr0LFDMfbXdP2l8vYkYqNyMVsls2wSRyLN+bB&#10;8RQ8IZv69HHzxLzrZy9ibD/b3ULpJ6qb1oNt8jR2toxWtjEpD7j2F2whSC/W3PN7tjrs9ekPAAAA&#10;//8DAFBLAwQUAAYACAAAACEAQ84ZOuAAAAAIAQAADwAAAGRycy9kb3ducmV2LnhtbEyPQU+DQBSE&#10;7yb+h80z8WaXtpEC8mgajYkHE2nVtMcFnkDKviXs0uK/d3uqx8lMZr5J15PuxIkG2xpGmM8CEMSl&#10;qVquEb4+Xx8iENYprlRnmBB+ycI6u71JVVKZM2/ptHO18CVsE4XQONcnUtqyIa3szPTE3vsxg1bO&#10;y6GW1aDOvlx3chEEodSqZb/QqJ6eGyqPu1Ej7Vr3RwmtLj52+DQEQX6zHMr0BkcFyN96OfnB4/h4&#10;7gAAAP//AwBQSwMEFAAGAAgAAAAhAEPOGTrgAAAACAEAAA8AAABkcnMvZG93bnJldi54bWxMj0FP&#10;g0EUhO8m/oftM/Fml7aRAvJomo2JBxNp1bTHBZ5Ayr4l7NLiv3d7qsfJTGa+SdeT7sSJBtsaRpjP&#10;AhDEpalarhG+Pl8fIhDWKa5UZ5gQfsnCOru9SVVSmTNv6bRztfAlbBOF0DjXJ1LasiGt7Mz0xN77&#10;MYNWzsuhltWgzr5cd3IRBKHUqmW/0KienhsqD7tRI+y3q+JljD/e86P5nod5fiC7eUO8v5s2TyAc&#10;Te4ahgu+R4fMMxVm5MqKDsEfcQirKIpBXOxwuQBRIMSPyxBklsr/B7I/AAAA//8DAFBLAQItABQA&#10;BgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1s&#10;UEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxz&#10;UEsBAi0AFAAGAAgAAAAhAD1GScPJAgAAHwYAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2Mu&#10;eG1sUEsBAi0AFAAGAAgAAAAhAEPOGTrgAAAACAEAAA8AAAAAAAAAAAAAAAAAIwUAAGRycy9kb3du&#10;cmV2LnhtbFBLBQYAAAAABAAEAPMAAAAwBgAAAAA=&#10;" filled="f" strokecolor="#70ad47 [3209]" strokeweight="3pt"><w10:wrap anchorx="margin"/></v:roundrect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251640822" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7D458480" wp14:editId="12205815">
'@

if ($xml.IndexOf($oldFirstRunStart) -lt 0) {
    throw "Anchor marker for the first picture run (relativeHeight=251640822) was not found."
}
$xml = $xml.Replace($oldFirstRunStart, $newFirstRunStart)

# --- 2) Stamp every other existing <wp:anchor> with its new anchorId/editId pair. ---
# relativeHeight=251667456
$old0 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251667456" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new0 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251667456" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="6428A26B" wp14:editId="05EBCB1A">
'@

if ($xml.IndexOf($old0) -lt 0) {
    throw "Anchor marker for relativeHeight=251667456 was not found."
}
$xml = $xml.Replace($old0, $new0)

# relativeHeight=251674624
$old1 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251674624" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="317A1989" wp14:editId="34327962">
'@

$new1 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251674624" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="41DB4AF5" wp14:editId="14D780C8">
'@

if ($xml.IndexOf($old1) -lt 0) {
    throw "Anchor marker for relativeHeight=251674624 was not found."
}
$xml = $xml.Replace($old1, $new1)

# relativeHeight=251676672
$old2 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251676672" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3455F6BD" wp14:editId="60E76DE0">
'@

$new2 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251676672" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="2A744216" wp14:editId="49182021">
'@

if ($xml.IndexOf($old2) -lt 0) {
    throw "Anchor marker for relativeHeight=251676672 was not found."
}
$xml = $xml.Replace($old2, $new2)

# relativeHeight=251677696
$old3 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251677696" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="254EC924" wp14:editId="5A8742A4">
'@

$new3 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251677696" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3D957F89" wp14:editId="2B4A73AD">
'@

if ($xml.IndexOf($old3) -lt 0) {
    throw "Anchor marker for relativeHeight=251677696 was not found."
}
$xml = $xml.Replace($old3, $new3)

# relativeHeight=251682816
$old4 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251682816" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7DAF7FCC" wp14:editId="50C6AB6C">
'@

$new4 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251682816" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="62448392" wp14:editId="79785C11">
'@

if ($xml.IndexOf($old4) -lt 0) {
    throw "Anchor marker for relativeHeight=251682816 was not found."
}
$xml = $xml.Replace($old4, $new4)

# relativeHeight=251681792
$old5 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251681792" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="14B59899" wp14:editId="15ABE8B8">
'@

$new5 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251681792" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1970E4D3" wp14:editId="1B0734D1">
'@

if ($xml.IndexOf($old5) -lt 0) {
    throw "Anchor marker for relativeHeight=251681792 was not found."
}
$xml = $xml.Replace($old5, $new5)

# relativeHeight=251680768
$old6 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251680768" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="292BEE7F" wp14:editId="15694E17">
'@

$new6 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251680768" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7E81F84F" wp14:editId="5B1FEC5A">
'@

if ($xml.IndexOf($old6) -lt 0) {
    throw "Anchor marker for relativeHeight=251680768 was not found."
}
$xml = $xml.Replace($old6, $new6)

# relativeHeight=251679744
$old7 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251679744" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3B28C6E7" wp14:editId="26636DC7">
'@

$new7 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251679744" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="72BB2A0F" wp14:editId="11497AD6">
'@

if ($xml.IndexOf($old7) -lt 0) {
    throw "Anchor marker for relativeHeight=251679744 was not found."
}
$xml = $xml.Replace($old7, $new7)

# relativeHeight=251656190
$old8 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251656190" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new8 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251656190" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4AFA211C" wp14:editId="78D4FFF1">
'@

if ($xml.IndexOf($old8) -lt 0) {
    throw "Anchor marker for relativeHeight=251656190 was not found."
}
$xml = $xml.Replace($old8, $new8)

# relativeHeight=251653115
$old9 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251653115" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="63D84C8A" wp14:editId="7BFD2B50">
'@

$new9 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251653115" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1E6BB96A" wp14:editId="25B4CA89">
'@

if ($xml.IndexOf($old9) -lt 0) {
    throw "Anchor marker for relativeHeight=251653115 was not found."
}
$xml = $xml.Replace($old9, $new9)

# relativeHeight=251654140
$old10 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251654140" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4191EA86" wp14:editId="15DEDCD1">
'@

$new10 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251654140" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="23CAB7FF" wp14:editId="33017D64">
'@

if ($xml.IndexOf($old10) -lt 0) {
    throw "Anchor marker for relativeHeight=251654140 was not found."
}
$xml = $xml.Replace($old10, $new10)

# relativeHeight=251652090
$old11 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251652090" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7DE029A9" wp14:editId="1BAB2B56">
'@

$new11 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251652090" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5B34314E" wp14:editId="7C545E37">
'@

if ($xml.IndexOf($old11) -lt 0) {
    throw "Anchor marker for relativeHeight=251652090 was not found."
}
$xml = $xml.Replace($old11, $new11)

# relativeHeight=251651065
$old12 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251651065" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5D57616B" wp14:editId="7C412CAB">
'@

$new12 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251651065" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4BAF5EDD" wp14:editId="21D719EE">
'@

if ($xml.IndexOf($old12) -lt 0) {
    throw "Anchor marker for relativeHeight=251651065 was not found."
}
$xml = $xml.Replace($old12, $new12)

# relativeHeight=251691008
$old13 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251691008" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="35BF5FEF" wp14:editId="09D01D2F">
'@

$new13 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251691008" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="540A8A8A" wp14:editId="3371AA1F">
'@

if ($xml.IndexOf($old13) -lt 0) {
    throw "Anchor marker for relativeHeight=251691008 was not found."
}
$xml = $xml.Replace($old13, $new13)

# relativeHeight=251692032
$old14 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251692032" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5B51E967" wp14:editId="7074034F">
'@

$new14 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251692032" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="2A9846DF" wp14:editId="2EE3AA29">
'@

if ($xml.IndexOf($old14) -lt 0) {
    throw "Anchor marker for relativeHeight=251692032 was not found."
}
$xml = $xml.Replace($old14, $new14)

# relativeHeight=251693056
$old15 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251693056" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5049940B" wp14:editId="31C15BBD">
'@

$new15 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251693056" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="6F03B5DE" wp14:editId="7C64F0B8">
'@

if ($xml.IndexOf($old15) -lt 0) {
    throw "Anchor marker for relativeHeight=251693056 was not found."
}
$xml = $xml.Replace($old15, $new15)

# relativeHeight=251694080
$old16 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251694080" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new16 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251694080" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4773B269" wp14:editId="3C9051D1">
'@

if ($xml.IndexOf($old16) -lt 0) {
    throw "Anchor marker for relativeHeight=251694080 was not found."
}
$xml = $xml.Replace($old16, $new16)

# relativeHeight=251699200
$old17 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251699200" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new17 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251699200" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="079F9971" wp14:editId="10CE47B7">
'@

if ($xml.IndexOf($old17) -lt 0) {
    throw "Anchor marker for relativeHeight=251699200 was not found."
}
$xml = $xml.Replace($old17, $new17)

# relativeHeight=251697152
$old18 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251697152" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new18 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251697152" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="54A1EBC8" wp14:editId="62DF2940">
'@

if ($xml.IndexOf($old18) -lt 0) {
    throw "Anchor marker for relativeHeight=251697152 was not found."
}
$xml = $xml.Replace($old18, $new18)

# relativeHeight=251696128
$old19 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251696128" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new19 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251696128" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="263E7509" wp14:editId="451CAC9E">
'@

if ($xml.IndexOf($old19) -lt 0) {
    throw "Anchor marker for relativeHeight=251696128 was not found."
}
$xml = $xml.Replace($old19, $new19)

# relativeHeight=251706368
$old20 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251706368" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="30973355" wp14:editId="339901FF">
'@

$new20 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251706368" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7816F4C6" wp14:editId="23B945CB">
'@

if ($xml.IndexOf($old20) -lt 0) {
    throw "Anchor marker for relativeHeight=251706368 was not found."
}
$xml = $xml.Replace($old20, $new20)

# relativeHeight=251708416
$old21 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251708416" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="791D93F1" wp14:editId="7EC18F9D">
'@

$new21 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251708416" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="15B74BF3" wp14:editId="33E95C88">
'@

if ($xml.IndexOf($old21) -lt 0) {
    throw "Anchor marker for relativeHeight=251708416 was not found."
}
$xml = $xml.Replace($old21, $new21)

# relativeHeight=251698176
$old22 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251698176" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new22 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251698176" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1A8F8A2D" wp14:editId="0FB49FBE">
'@

if ($xml.IndexOf($old22) -lt 0) {
    throw "Anchor marker for relativeHeight=251698176 was not found."
}
$xml = $xml.Replace($old22, $new22)

# relativeHeight=251710464
$old23 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251710464" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3429F1AB" wp14:editId="4828F461">
'@

$new23 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251710464" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1F60E2C4" wp14:editId="62E82928">
'@

if ($xml.IndexOf($old23) -lt 0) {
    throw "Anchor marker for relativeHeight=251710464 was not found."
}
$xml = $xml.Replace($old23, $new23)

# relativeHeight=251641847
$old24 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251641847" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new24 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251641847" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7F7FC6F5" wp14:editId="6400A6A1">
'@

if ($xml.IndexOf($old24) -lt 0) {
    throw "Anchor marker for relativeHeight=251641847 was not found."
}
$xml = $xml.Replace($old24, $new24)

# relativeHeight=251646967
$old25 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251646967" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="11CE4D02" wp14:editId="326D135A">
'@

$new25 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251646967" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="0805519F" wp14:editId="4F190027">
'@

if ($xml.IndexOf($old25) -lt 0) {
    throw "Anchor marker for relativeHeight=251646967 was not found."
}
$xml = $xml.Replace($old25, $new25)

# relativeHeight=251642871
$old26 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251642871" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new26 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251642871" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="0D27A546" wp14:editId="44FC1E62">
'@

if ($xml.IndexOf($old26) -lt 0) {
    throw "Anchor marker for relativeHeight=251642871 was not found."
}
$xml = $xml.Replace($old26, $new26)

# relativeHeight=251647991
$old27 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251647991" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="285A8C58" wp14:editId="34C6C271">
'@

$new27 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251647991" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="17267909" wp14:editId="39FE032B">
'@

if ($xml.IndexOf($old27) -lt 0) {
    throw "Anchor marker for relativeHeight=251647991 was not found."
}
$xml = $xml.Replace($old27, $new27)

# relativeHeight=251643895
$old28 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251643895" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new28 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251643895" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="60CC6B9A" wp14:editId="1D48CAC4">
'@

if ($xml.IndexOf($old28) -lt 0) {
    throw "Anchor marker for relativeHeight=251643895 was not found."
}
$xml = $xml.Replace($old28, $new28)

# relativeHeight=251649015
$old29 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251649015" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="491651CD" wp14:editId="784D03BB">
'@

$new29 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251649015" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="0C977030" wp14:editId="28B790A5">
'@

if ($xml.IndexOf($old29) -lt 0) {
    throw "Anchor marker for relativeHeight=251649015 was not found."
}
$xml = $xml.Replace($old29, $new29)

# relativeHeight=251644919
$old30 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251644919" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new30 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251644919" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="2534CE0E" wp14:editId="75D1A6D4">
'@

if ($xml.IndexOf($old30) -lt 0) {
    throw "Anchor marker for relativeHeight=251644919 was not found."
}
$xml = $xml.Replace($old30, $new30)

# relativeHeight=251645943
$old31 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251645943" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1">
'@

$new31 = @'
<wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251645943" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="06FA189E" wp14:editId="438F13C7">
'@

if ($xml.IndexOf($old31) -lt 0) {
    throw "Anchor marker for relativeHeight=251645943 was not found."
}
$xml = $xml.Replace($old31, $new31)

# --- 3) Write the modified package XML back to the document. ---
$rng.WordOpenXML = $xml

Write-Output "Applied Chess design.docx box-fix edit successfully."
